$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix an existing row's academic year (transcript correction)
$ws.Cells.Item(21, 1).Value = 2015

# Row 52 is filled in first (matches authoring order reflected in shared strings),
# capturing a Zoo699 transfer credit from 2010
$ws.Cells.Item(52, 1).Value = 2010
$ws.Cells.Item(52, 2).Value = "TRANSFER"
$ws.Cells.Item(52, 3).Value = "Zoo699"
$ws.Cells.Item(52, 4).Value = 1
$ws.Cells.Item(52, 7).Value = "S"

# Existing transfer rows 49 and 50 get reclassified with actual course names
# and updated credit amounts, plus a new row 51 for a third AP transfer course
$ws.Cells.Item(49, 1).Value = 2012
$ws.Cells.Item(49, 2).Value = "TRANSFER"
$ws.Cells.Item(49, 3).Value = "APBio"
$ws.Cells.Item(49, 4).Value = 3
$ws.Cells.Item(49, 7).Value = "S"

$ws.Cells.Item(50, 1).Value = 2012
$ws.Cells.Item(50, 2).Value = "TRANSFER"
$ws.Cells.Item(50, 3).Value = "APEngLit"
$ws.Cells.Item(50, 4).Value = 3
$ws.Cells.Item(50, 7).Value = "S"

$ws.Cells.Item(51, 1).Value = 2012
$ws.Cells.Item(51, 2).Value = "TRANSFER"
$ws.Cells.Item(51, 3).Value = "APEnvSci"
$ws.Cells.Item(51, 4).Value = 3
$ws.Cells.Item(51, 7).Value = "S"

# Update the view state to reflect scrolling down to the newly added rows
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E49").Select()
